$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New fund flow data rows to append (asOf date as text, value as number)
$data = @(
    @("2023-10-30", 0.45571336),
    @("2023-10-31", 50.4175),
    @("2023-11-01", 5.76),
    @("2023-11-02", 38.8395),
    @("2023-11-03", 46.016)
)

$startRow = 210
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i

    $cellA = $ws.Cells.Item($r, 1)
    # Force the date-looking string to be stored as literal text (not parsed
    # into a date serial number), matching the existing asOf column, then
    # reset the style back to the workbook default so no extra style index
    # is attached to the cell.
    $cellA.NumberFormat = "@"
    $cellA.Value = $data[$i][0]
    $cellA.Style = "Normal"

    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = $data[$i][1]
}
